$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 512
$ws1.Range("F5").Value = 912
$ws1.Range("F6").Value = 148
$ws1.Range("F7").Value = 917
$ws1.Range("F8").Value = 714
$ws1.Range("F9").Value = 168
$ws1.Range("F11").Value = 77
$ws1.Range("F13").Value = 248
$ws1.Range("F14").Value = 540
$ws1.Range("F15").Value = 483
$ws1.Range("F16").Value = 1275
$ws1.Range("F17").Value = 111
$ws1.Range("F18").Value = 414
$ws1.Range("F19").Value = 1035
$ws1.Range("F20").Value = 2767
$ws1.Range("F21").Value = 1245
$ws1.Range("F22").Value = 635
$ws1.Range("F24").Value = 1232
$ws1.Range("F26").Value = 957
$ws1.Range("F27").Value = 313
$ws1.Range("F28").Value = 398
$ws1.Range("F29").Value = 1292

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 715

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 715
$ws4.Range("F3").Value = 512
$ws4.Range("F12").Value = 912
$ws4.Range("F13").Value = 148
$ws4.Range("F15").Value = 917
$ws4.Range("F16").Value = 714
$ws4.Range("F17").Value = 168
$ws4.Range("F23").Value = 77
$ws4.Range("F26").Value = 248
$ws4.Range("F27").Value = 540
$ws4.Range("F28").Value = 483
$ws4.Range("F29").Value = 1275
$ws4.Range("F30").Value = 111
$ws4.Range("F31").Value = 414
$ws4.Range("F32").Value = 1035
$ws4.Range("F33").Value = 2767
$ws4.Range("F34").Value = 1245
$ws4.Range("F35").Value = 635
$ws4.Range("F37").Value = 1232
$ws4.Range("F40").Value = 957
$ws4.Range("F41").Value = 313
$ws4.Range("F42").Value = 398
$ws4.Range("F43").Value = 1292
